$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose values look numeric so Excel
# keeps them as strings (matching the original inlineStr data),
# not auto-converted numbers.

$ws.Range('D2').Value = '63.570.87'
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = '3.089.42'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.26'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.34'
$ws.Range('E6').Value = '  +3.62%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.550'
$ws.Range('E8').Value = '  +3.46%  '
$ws.Range('D9').Value = '3.080.46'
$ws.Range('E9').Value = '  -2.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.157'
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.86'
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.51'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('D15').Value = '3.600.81'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '63.568.88'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.14'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').Value = '3.084.77'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.79'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.65'
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.12'
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.82'
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.70'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('E29').Value = '  -2.16%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  -2.92%  '
$ws.Range('E32').Value = '  -3.81%  '
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.43'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '0.0₃0849'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.36'
$ws.Range('E37').Value = '  +2.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.10'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.22'
$ws.Range('E39').Value = '  -5.36%  '
$ws.Range('E40').Value = '  +0.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.64'
$ws.Range('E41').Value = '  -2.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '442.92'
$ws.Range('E42').Value = '  -3.44%  '
$ws.Range('E43').Value = '  -4.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0364'
$ws.Range('E44').Value = '  -2.73%  '
$ws.Range('D45').Value = '2.820.10'
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.11'
$ws.Range('E46').Value = '  -1.69%  '
$ws.Range('B47').Value = 'Kaspa'
$ws.Range('C47').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.109'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.20'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.15'
$ws.Range('E50').Value = '  +2.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.26'
$ws.Range('E51').Value = '  +0.41%  '
